$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column H: mirror of column A (distance values), with the same
#     bold/bordered style used by A/D --------------------------------------
$ws.Range("A2:A21").Copy()
$ws.Range("H2:H21").PasteSpecial(-4104)  # xlPasteAll (values)
$ws.Range("A2:A21").Copy()
$ws.Range("H2:H21").PasteSpecial(-4122)  # xlPasteFormats (style s="1")
$excel.CutCopyMode = 0

# --- New column I: distance recovered from the measured dB value ----------
#     POWER(10,(E{row}-$F$1)/(10*$G$1))*100
$ws.Range("I2").Formula = '=POWER(10,(E2-$F$1)/(10*$G$1))*100'
$ws.Range("I3:I21").Formula = '=POWER(10,(E3-$F$1)/(10*$G$1))*100'

# --- Chart was dragged/resized to a new spot on the sheet -----------------
$co = $ws.ChartObjects().Item(1)
$co.Left = 947.3
$co.Top = 60.15
$co.Width = 608.675
$co.Height = 316.8

# --- Selection moved to M9 as in the edited workbook -----------------------
$null = $ws.Range("M9").Select()
